# "Tried fixing length and height calculation"
#
# Sheet "1_filtered_data":
#   - Column I ("length", C-E or E-C) gets an extra +/-76 offset baked into
#     the formula (instead of the separate constant living implicitly in
#     the raw pixel values).
#   - Column K ("center") is rewritten from a straight midpoint
#     ((C+E)/2) to a midpoint measured from the 1614 reference point
#     (1614-(C+E)/2).
#   - The active selection moves from H2 to M32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1_filtered_data")

# Clear the existing formulas first so the sheet's shared-formula
# bookkeeping for columns I/K is rebuilt from scratch rather than being
# stitched onto the old si="0"/si="2"/si="4"/si="5" groups.
$ws.Range("I2:I30").ClearContents()
$ws.Range("K2:K30").ClearContents()

# --- Column I ------------------------------------------------------------
# Row 2: standalone formula (own <f>, not shared) - matches the original
# layout where I2 was never part of the I3:I20 shared group.
$ws.Range("I2").Formula = "=C2-E2-76"

# Rows 3-20 share one formula (C-E-76).
$ws.Range("I3:I20").Formula = "=C3-E3-76"

# Row 21: standalone formula, opposite sign convention (E-C+76).
$ws.Range("I21").Formula = "=E21-C21+76"

# Rows 22-30 share the opposite-sign formula.
$ws.Range("I22:I30").Formula = "=E22-C22+76"

# --- Column K --------------------------------------------------------------
# Row 2: standalone formula.
$ws.Range("K2").Formula = "=1614-(C2+E2)/2"

# Rows 3-20 share one formula.
$ws.Range("K3:K20").Formula = "=1614-(C3+E3)/2"

# Row 21: standalone formula (same shape, kept separate like the source).
$ws.Range("K21").Formula = "=1614-(C21+E21)/2"

# Rows 22-30 share one formula.
$ws.Range("K22:K30").Formula = "=1614-(C22+E22)/2"

# --- Selection -------------------------------------------------------------
$ws.Activate()
$ws.Range("M32").Select()
